$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the "Contact Information" Heading2 paragraph entirely.
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Contact Information") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Helper: split a run of text into three runs, wrapping the middle `word`
# with <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/>,
# exactly mirroring the OOXML produced by Word's background grammar check.
# ---------------------------------------------------------------------------
function Split-WithProofErr($fullText, $word_) {
    $d = $script:d
    $rng = $d.Content
    $found = $rng.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $fullText
        return
    }

    $start = $rng.Start
    $end = $rng.End
    $target = $d.Range($start, $end)

    $wordIdx = $fullText.IndexOf($word_)
    $before = $fullText.Substring(0, $wordIdx)
    $after = $fullText.Substring($wordIdx + $word_.Length)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body><w:p>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve">' + $before + '</w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t>' + $word_ + '</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve">' + $after + '</w:t></w:r>' +
      '</w:p></w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

$script:d = $d

# ---------------------------------------------------------------------------
# 2) "Deployed URL shortener for free and personal use amongst friends and
#    family" -> split around "shortener"
# ---------------------------------------------------------------------------
Split-WithProofErr "Deployed URL shortener for free and personal use amongst friends and family" "shortener"

# ---------------------------------------------------------------------------
# 3) "Researched accessibility gaps ... support alt text and visual
#    descriptions" -> split around "alt"
# ---------------------------------------------------------------------------
Split-WithProofErr "Researched accessibility gaps in the IIIF spec and common image viewers, drafting proposed metadata enhancements to support alt text and visual descriptions" "alt"

Write-Host "Done"
